# Add 2022-Q4 data
#
# 1) Insert a brand-new worksheet named "2022-Q4" right before "2022-Q3"
#    (i.e. right after "总计"), and fill it with the new quarterly holdings
#    table (mirrors the other "2022-QN" sheets: bold/centered/bordered
#    header row + bold/centered/bordered index column).
# 2) Update the "总计" (totals) summary sheet: shift the existing rows
#    down by one and insert the new 2022-Q4 totals at the top (row 2),
#    re-adding the 2022-Q1 row that falls off the bottom at row 5.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the "2022-Q4" worksheet, positioned before "2022-Q3"
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$wb.Worksheets.Add($q3) | Out-Null
$wb.Worksheets.Item(2).Name = "2022-Q4"

$headers = @("基金代码","基金名称","基金规模","股票总仓位","仓位占比","持有市值(亿元)","仓位排名")

$rows = @(
    @("009693","富国积极成长一年定期开放混合","12.47","98.05","5.90","0.7357",2),
    @("014273","广发北交所精选两年定开混合A","3.23","83.79","6.50","0.2100",3),
    @("014269","嘉实北交所精选两年定期混合A","2.65","94.48","6.64","0.1760",2),
    @("014275","易方达北交所精选两年定开混合A","3.47","70.57","3.12","0.1083",5),
    @("014663","富国创新发展两年定期开放混合A","2.26","82.02","3.51","0.0793",3),
    @("014274","广发北交所精选两年定开混合C","0.81","83.79","6.50","0.0526",3),
    @("014270","嘉实北交所精选两年定期混合C","0.52","94.48","6.64","0.0345",2),
    @("014276","易方达北交所精选两年定开混合C","0.90","70.57","3.12","0.0281",5),
    @("014664","富国创新发展两年定期开放混合C","0.33","82.02","3.51","0.0116",3)
)
$rowCount = $rows.Length

# Seed every cell we're about to touch so the sheet has real rows/cells to
# copy formatting onto (re-fetch the sheet by name each time rather than
# re-using an old reference, since stale worksheet handles silently drop
# pasted formatting in this host).
$wb.Worksheets.Item("2022-Q4").Range("B1:H1").Value = "seed" | Out-Null
$wb.Worksheets.Item("2022-Q4").Range("A2:A" + ($rowCount + 1)).Value = 0 | Out-Null

# Mirror the bold/centered/thin-bordered look used by the header row and
# the "A" index column on the sibling "2022-QN" sheets by copying their
# cell formatting (values untouched) onto the new sheet.
$wb.Worksheets.Item("2022-Q3").Range("B1:H1").Copy()
$wb.Worksheets.Item("2022-Q4").Range("B1:H1").PasteSpecial(-4122)

$wb.Worksheets.Item("2022-Q3").Range("A2").Copy()
$wb.Worksheets.Item("2022-Q4").Range(("A2:A" + ($rowCount + 1))).PasteSpecial(-4122)

# Header row text
for ($c = 0; $c -lt $headers.Length; $c++) {
    $wb.Worksheets.Item("2022-Q4").Cells.Item(1, $c + 2).Value = $headers[$c]
}

# Data rows
for ($r = 0; $r -lt $rowCount; $r++) {
    $rowIndex = $r + 2
    $values = $rows[$r]
    $sheet = $wb.Worksheets.Item("2022-Q4")

    $sheet.Cells.Item($rowIndex, 1).Value = $r

    # Columns B-G are text in the source data (fund codes with leading
    # zeros, and numeric-looking strings whose trailing zeros / decimal
    # formatting must be preserved verbatim) - force text with a leading
    # apostrophe so Excel doesn't re-interpret them as numbers.
    $sheet.Cells.Item($rowIndex, 2).Value = "'" + $values[0]
    $sheet.Cells.Item($rowIndex, 3).Value = $values[1]
    $sheet.Cells.Item($rowIndex, 4).Value = "'" + $values[2]
    $sheet.Cells.Item($rowIndex, 5).Value = "'" + $values[3]
    $sheet.Cells.Item($rowIndex, 6).Value = "'" + $values[4]
    $sheet.Cells.Item($rowIndex, 7).Value = "'" + $values[5]
    $sheet.Cells.Item($rowIndex, 8).Value = $values[6]
}

# ---------------------------------------------------------------------
# 2. Update the "总计" summary sheet with the new quarter on top
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Row 5 is brand new - give it the same formatting as the existing index
# rows before filling in the (shifted-down) 2022-Q1 totals.
$total.Range("A4").Copy()
$total.Range("A5").PasteSpecial(-4122)

# Existing rows 2-4 (Q3, Q2, Q1) shift down to rows 3-5; insert the new
# Q4 totals at row 2.
$total.Cells.Item(5, 1).Value = 3
$total.Cells.Item(5, 2).Value = "2022-Q1"
$total.Cells.Item(5, 3).Value = 5
$total.Cells.Item(5, 4).Value = 1.01

$total.Cells.Item(4, 1).Value = 2
$total.Cells.Item(4, 2).Value = "2022-Q2"
$total.Cells.Item(4, 3).Value = 5
$total.Cells.Item(4, 4).Value = 0.91

$total.Cells.Item(3, 1).Value = 1
$total.Cells.Item(3, 2).Value = "2022-Q3"
$total.Cells.Item(3, 3).Value = 7
$total.Cells.Item(3, 4).Value = 0.93

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q4"
$total.Cells.Item(2, 3).Value = 9
$total.Cells.Item(2, 4).Value = 1.44
